$wb = $excel.ActiveWorkbook

# The "Neighbors" sheet (2nd sheet) had its row 2 ("Obstructed (T/F)" row, all FALSE)
# deleted entirely; subsequent rows shift up by one.
$wsNeighbors = $wb.Worksheets.Item("Neighbors")
$wsNeighbors.Rows.Item(2).Delete() | Out-Null

# Update selection on the Neighbors sheet to match the new layout.
$wsNeighbors.Application.ActiveWindow.ScrollColumn = 1
$wsNeighbors.Range("A2:I13").Select() | Out-Null

# Make "Neighbors" the active (selected) sheet/tab.
$wsNeighbors.Activate()

$wb.Save()
